# Apply updated crypto price/volume data (GitHub Actions refresh)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "56.911.86"
$ws.Range("E2").Value = "  +0.14%  "

$ws.Range("D3").Value = "2.418.98"
$ws.Range("E3").Value = "  -1.47%  "

$ws.Range("E4").Value = "  +0.25%  "

$ws.Range("D5").Value = "'487.49"

$ws.Range("D6").Value = "'153.22"
$ws.Range("E6").Value = "  +0.90%  "

$ws.Range("B7").Value = "USDC"
$ws.Range("C7").Value = "https://coinranking.com/coin/aKzUVe4Hh_CON+usdc-usdc"
$ws.Range("D7").Value = "'0.997"
$ws.Range("E7").Value = "  -0.03%  "

$ws.Range("B8").Value = "XRP"
$ws.Range("C8").Value = "https://coinranking.com/coin/-l8Mn2pVlRs-p+xrp-xrp"
$ws.Range("D8").Value = "'0.608"
$ws.Range("E8").Value = "  +18.26%  "

$ws.Range("D9").Value = "2.443.38"
$ws.Range("E9").Value = "  -0.68%  "

$ws.Range("D10").Value = "'6.16"
$ws.Range("E10").Value = "  +8.39%  "

$ws.Range("D11").Value = "'0.0995"
$ws.Range("E11").Value = "  -0.01%  "

$ws.Range("D12").Value = "'0.331"
$ws.Range("E12").Value = "  -0.67%  "

$ws.Range("E13").Value = "  +1.20%  "

$ws.Range("D14").Value = "2.841.48"
$ws.Range("E14").Value = "  -2.04%  "

$ws.Range("D15").Value = "56.979.17"
$ws.Range("E15").Value = "  -0.25%  "

$ws.Range("D16").Value = "'20.49"
$ws.Range("E16").Value = "  -2.13%  "

$ws.Range("D17").Value = "'0.0000132"
$ws.Range("E17").Value = "  -3.29%  "

$ws.Range("D18").Value = "2.460.44"
$ws.Range("E18").Value = "  -0.21%  "

$ws.Range("D19").Value = "'4.58"
$ws.Range("E19").Value = "  +0.86%  "

$ws.Range("D20").Value = "'322.24"
$ws.Range("E20").Value = "  -0.23%  "

$ws.Range("D21").Value = "'9.99"
$ws.Range("E21").Value = "  -1.09%  "

$ws.Range("D22").Value = "'0.997"
$ws.Range("E22").Value = "  +0.08%  "

$ws.Range("D23").Value = "'5.90"
$ws.Range("E23").Value = "  +1.31%  "

$ws.Range("D24").Value = "'57.74"
$ws.Range("E24").Value = "  -0.32%  "

$ws.Range("D25").Value = "'0.997"
$ws.Range("E25").Value = "  -0.31%  "

$ws.Range("D26").Value = "'0.399"
$ws.Range("E26").Value = "  -1.43%  "

$ws.Range("D27").Value = "'0.158"
$ws.Range("E27").Value = "  -2.85%  "

$ws.Range("D28").Value = "2.534.60"
$ws.Range("E28").Value = "  -1.59%  "

$ws.Range("D29").Value = "'7.24"
$ws.Range("E29").Value = "  -4.02%  "

$ws.Range("D30").Value = "0.0₃0784"
$ws.Range("E30").Value = "  -2.37%  "

$ws.Range("E31").Value = "  -0.01%  "

$ws.Range("B32").Value = "Monero"
$ws.Range("C32").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D32").Value = "'150.56"
$ws.Range("E32").Value = "  +0.06%  "

$ws.Range("B33").Value = "EthereumClassic"
$ws.Range("C33").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D33").Value = "'18.62"
$ws.Range("E33").Value = "  +1.98%  "

$ws.Range("D34").Value = "'1.51"
$ws.Range("E34").Value = "  +0.03%  "

$ws.Range("D35").Value = "'5.27"
$ws.Range("E35").Value = "  +1.34%  "

$ws.Range("D36").Value = "'3.76"
$ws.Range("E36").Value = "  -0.03%  "

$ws.Range("D37").Value = "'1.12"
$ws.Range("E37").Value = "  -1.72%  "

$ws.Range("D38").Value = "'0.811"
$ws.Range("E38").Value = "  -8.39%  "

$ws.Range("B39").Value = "OKB"
$ws.Range("C39").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D39").Value = "'34.01"
$ws.Range("E39").Value = "  -0.63%  "

$ws.Range("D40").Value = "'0.100"
$ws.Range("E40").Value = "  +5.19%  "

$ws.Range("B41").Value = "Filecoin"
$ws.Range("C41").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D41").Value = "'3.50"
$ws.Range("E41").Value = "  +0.36%  "

$ws.Range("D42").Value = "'1.36"
$ws.Range("E42").Value = "  -2.16%  "

$ws.Range("B43").Value = "Bittensor"
$ws.Range("C43").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D43").Value = "'279.73"
$ws.Range("E43").Value = "  +7.17%  "

$ws.Range("D44").Value = "'0.995"
$ws.Range("E44").Value = "  -0.16%  "

$ws.Range("D45").Value = "'0.600"
$ws.Range("E45").Value = "  -0.99%  "

$ws.Range("D46").Value = "'0.0529"
$ws.Range("E46").Value = "  -4.62%  "

$ws.Range("D47").Value = "'10.20"
$ws.Range("E47").Value = "  -0.44%  "

$ws.Range("E48").Value = "  -0.16%  "

$ws.Range("D49").Value = "'4.51"
$ws.Range("E49").Value = "  -6.45%  "

$ws.Range("B50").Value = "EnergySwap"
$ws.Range("C50").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D50").Value = "'17.76"
$ws.Range("E50").Value = "  +0.07%  "

$ws.Range("B51").Value = "Maker"
$ws.Range("C51").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D51").Value = "1.897.86"
$ws.Range("E51").Value = "  +3.20%  "
